# Apply the Alvearie FHIR IG deployment update to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Metadata" worksheet updates
# ---------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date update
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Replace the first "Contact" row with a new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the now-duplicate second "Contact" row entirely, shifting
# everything below it up by one row.
$meta.Rows.Item(11).Delete()

# ---------------------------------------------------------------
# 2. "Elements" worksheet updates
# ---------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) gets a real Short / Definition
# instead of the generic placeholder text.
$elements.Range("K2").Value = "Source Record Type"
$elements.Range("L2").Value = "Either the data model type or schema type that generates this FHIR resource"
